# Upgrade left table until javakheti
# - Rename the worksheet from "1" to "Martvili"
# - Remove the blank row 8 (Note row moves from row 9 up to row 8)
# - Mark the "Urban" row (row 6, years 2010-2020) as confidential ("...")
# - Mark a handful of "Rural" row (row 7) cells as confidential/unavailable

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Martvili"

# Delete the empty row 8 so the Note row shifts from row 9 to row 8
$ws.Rows.Item(8).Delete()

# Urban row (row 6): years 2010-2020 (columns B-L) become confidential "..."
$ws.Range("B6:L6").Value = "..."

# Rural row (row 7): specific years become unavailable/confidential
$ws.Range("C7").Value = "..."
$ws.Range("F7:G7").Value = "..."
$ws.Range("J7:K7").Value = "…"
